$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A58").Value = "Daniel Azzolini"
$ws.Range("B58").Value = "Stefano Tita | Clitoriders"
$ws.Range("C58").Value = "Marco  Sartorelli | Modium"
$ws.Range("D58").Value = "Giacomo Gasparini | MAI UNA GIOIA"
$ws.Range("E58").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("F58").Value = "Davide  Bazzano  | iMontagna"
